$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Human"/"Ecosystem" categories (columns H = Type_bar2 and
# L = Type_CA) to the new "Social justice"/"Environmental justice"
# wording used for the panel legend added on the figures.
$map = @{
    "Human "    = "Social justice "
    "Ecosystem" = "Environmental justice"
    "Human"     = "Social justice"
}

# Seed the three replacement strings in the same order they end up in
# the workbook's shared-string table after the edit (preserve-space
# variant, then the bare "Social justice", then "Environmental
# justice") before sweeping the rest of the column, so new values line
# up cell-for-cell with the target workbook.
$ws.Range("H4").Value = "Social justice "
$ws.Range("H47").Value = "Social justice"
$ws.Range("H5").Value = "Environmental justice"

foreach ($col in @("H", "L")) {
    for ($row = 2; $row -le 56; $row++) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($null -ne $val -and $map.ContainsKey($val)) {
            $cell.Value = $map[$val]
        }
    }
}

# Update the saved selection/view state to match the author's final
# cursor position (H9, no scrolled-away top-left cell).
[void]$ws.Range("H9").Select()
